$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.035.23"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "3.862.74"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "700.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.45"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "3.861.63"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  -4.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.32"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "4.511.84"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "3.852.66"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "71.067.14"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "501.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.181"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "3.818.71"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.36%  "
$ws.Range("E40").Value = "  +8.59%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000314"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "416.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.50%  "
